$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure the cells we are about to write remain stored as literal text
# (matching original inlineStr cells), not auto-converted to numbers/percentages.
$textCells = @(
    "D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5", "E5", "G5", "E6", "G6", "D7", "E7",
    "G7", "D8", "E8", "G8", "D9", "E9", "G9", "D10", "E10", "G10", "D11", "E11", "G11", "D12",
    "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14", "D15", "E15", "G15", "D16", "E16",
    "G16", "D17", "E17", "G17", "D18", "E18", "G18", "G19", "D20", "E20", "G20", "E21", "G21",
    "D22", "E22", "G22", "D23", "E23", "G23", "D24", "E24", "G24", "D25", "E25", "G25", "D26",
    "E26", "G26", "D27", "E27", "G27", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35",
    "G36", "G37", "G38", "D39", "E39", "G39", "D40", "E40", "G40", "D41", "E41", "G41", "D42",
    "E42", "G42", "D43", "E43", "G43", "D44", "E44", "G44", "D45", "E45", "G45", "D46", "E46",
    "G46", "E47", "G47", "D48", "E48", "G48", "D49", "E49", "G49", "E50", "G50", "D51", "E51", "G51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "331.44"
$ws.Range("E2").Value = "0.35%"
$ws.Range("G2").Value = "3"
$ws.Range("D3").Value = "44.22"
$ws.Range("E3").Value = "7.40%"
$ws.Range("G3").Value = "3"
$ws.Range("D4").Value = "5.762"
$ws.Range("E4").Value = "1.13%"
$ws.Range("G4").Value = "3"
$ws.Range("D5").Value = "0.08302"
$ws.Range("E5").Value = "1.59%"
$ws.Range("G5").Value = "3"
$ws.Range("E6").Value = "0.86%"
$ws.Range("G6").Value = "3"
$ws.Range("D7").Value = "4.508"
$ws.Range("E7").Value = "-0.78%"
$ws.Range("G7").Value = "3"
$ws.Range("D8").Value = "1.965"
$ws.Range("E8").Value = "-5.20%"
$ws.Range("G8").Value = "3"
$ws.Range("D9").Value = "2.927"
$ws.Range("E9").Value = "-0.30%"
$ws.Range("G9").Value = "3"
$ws.Range("D10").Value = "0.9236"
$ws.Range("E10").Value = "-0.14%"
$ws.Range("G10").Value = "3"
$ws.Range("D11").Value = "0.1241"
$ws.Range("E11").Value = "-0.85%"
$ws.Range("G11").Value = "3"
$ws.Range("D12").Value = "0.1944"
$ws.Range("E12").Value = "-0.43%"
$ws.Range("G12").Value = "3"
$ws.Range("D13").Value = "0.09437"
$ws.Range("E13").Value = "1.65%"
$ws.Range("G13").Value = "3"
$ws.Range("D14").Value = "0.03971"
$ws.Range("E14").Value = "7.24%"
$ws.Range("G14").Value = "3"
$ws.Range("D15").Value = "0.1066"
$ws.Range("E15").Value = "1.05%"
$ws.Range("G15").Value = "3"
$ws.Range("D16").Value = "0.001304"
$ws.Range("E16").Value = "0.30%"
$ws.Range("G16").Value = "3"
$ws.Range("D17").Value = "0.006013"
$ws.Range("E17").Value = "-2.83%"
$ws.Range("G17").Value = "3"
$ws.Range("D18").Value = "3.535"
$ws.Range("E18").Value = "3.62%"
$ws.Range("G18").Value = "3"
$ws.Range("G19").Value = "3"
$ws.Range("D20").Value = "9.142"
$ws.Range("E20").Value = "9.92%"
$ws.Range("G20").Value = "3"
$ws.Range("E21").Value = "-0.47%"
$ws.Range("G21").Value = "3"
$ws.Range("D22").Value = "0.2579"
$ws.Range("E22").Value = "-2.75%"
$ws.Range("G22").Value = "3"
$ws.Range("D23").Value = "0.04428"
$ws.Range("E23").Value = "-0.43%"
$ws.Range("G23").Value = "3"
$ws.Range("D24").Value = "0.001244"
$ws.Range("E24").Value = "-2.29%"
$ws.Range("G24").Value = "3"
$ws.Range("D25").Value = "0.004413"
$ws.Range("E25").Value = "2.79%"
$ws.Range("G25").Value = "3"
$ws.Range("D26").Value = "0.0001194"
$ws.Range("E26").Value = "1.01%"
$ws.Range("G26").Value = "3"
$ws.Range("D27").Value = "0.0004002"
$ws.Range("E27").Value = "0.23%"
$ws.Range("G27").Value = "3"
$ws.Range("G28").Value = "3"
$ws.Range("G29").Value = "3"
$ws.Range("G30").Value = "3"
$ws.Range("G31").Value = "3"
$ws.Range("G32").Value = "3"
$ws.Range("G33").Value = "3"
$ws.Range("G34").Value = "3"
$ws.Range("G35").Value = "3"
$ws.Range("G36").Value = "3"
$ws.Range("G37").Value = "3"
$ws.Range("G38").Value = "3"
$ws.Range("D39").Value = "0.02829"
$ws.Range("E39").Value = "0.71%"
$ws.Range("G39").Value = "3"
$ws.Range("D40").Value = "0.05612"
$ws.Range("E40").Value = "2.39%"
$ws.Range("G40").Value = "3"
$ws.Range("D41").Value = "0.007940"
$ws.Range("E41").Value = "3.45%"
$ws.Range("G41").Value = "3"
$ws.Range("D42").Value = "0.1424"
$ws.Range("E42").Value = "0.44%"
$ws.Range("G42").Value = "3"
$ws.Range("D43").Value = "0.009077"
$ws.Range("E43").Value = "-4.03%"
$ws.Range("G43").Value = "3"
$ws.Range("D44").Value = "0.002126"
$ws.Range("E44").Value = "-0.31%"
$ws.Range("G44").Value = "3"
$ws.Range("D45").Value = "0.009934"
$ws.Range("E45").Value = "-17.00%"
$ws.Range("G45").Value = "3"
$ws.Range("D46").Value = "0.00007369"
$ws.Range("E46").Value = "7.15%"
$ws.Range("G46").Value = "3"
$ws.Range("E47").Value = "0.15%"
$ws.Range("G47").Value = "3"
$ws.Range("D48").Value = "0.003606"
$ws.Range("E48").Value = "11.58%"
$ws.Range("G48").Value = "3"
$ws.Range("D49").Value = "0.002285"
$ws.Range("E49").Value = "0.12%"
$ws.Range("G49").Value = "3"
$ws.Range("E50").Value = "0.15%"
$ws.Range("G50").Value = "3"
$ws.Range("D51").Value = "0.0002006"
$ws.Range("E51").Value = "0.15%"
$ws.Range("G51").Value = "3"

# Restore default cell style so no stray formatting/style index is introduced
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
